# Auto-generated edit script applying the Typhon_Profits.xlsx market-data refresh
# across all 8 leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$updates = @{
    "H64" = 4434.7144
    "I64" = 4043
    "K64" = 4043
    "M64" = -3795
    "H67" = 4434.7144
    "I67" = 4043
    "K67" = 4043
    "M67" = -3185
    "H76" = 2927065
    "I76" = 3167.6155
    "J76" = 9262176
    "K76" = 3167.6155
    "L76" = 9262176
    "M76" = -2852.6155
    "N76" = -9262806
    "H79" = 2927065
    "I79" = 3167.6155
    "J79" = 9262176
    "K79" = 3167.6155
    "L79" = 9262176
    "M79" = -2075.6155
    "N79" = -9264360
    "H100" = 3401
    "I100" = 3001.25
    "J100" = 5000
    "K100" = 3001.25
    "L100" = 5000
    "M100" = -2460.25
    "N100" = -6082
    "H116" = 6155.091
    "I116" = 3666.6667
    "J116" = 7088.25
    "K116" = 3666.6667
    "L116" = 7088.25
    "M116" = -224.6667000000002
    "N116" = -13972.25
    "H127" = 1214.6428
    "I127" = 769.625
    "K127" = 2308.875
    "M127" = 2651.125
    "H137" = 1614.8966
    "I137" = 1449.0526
    "J137" = 1930
    "K137" = 4347.1578
    "L137" = 5790
    "M137" = -1797.1578
    "N137" = -10890
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("ARM")

$updates = @{
    "H32" = 2812.7273
    "I32" = 2544.973
    "J32" = 4228
    "K32" = 2544.973
    "L32" = 4228
    "M32" = -2257.973
    "N32" = -4802
    "H122" = 3140
    "J122" = 2600
    "L122" = 7800
    "N122" = -12700
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("BSM")

$updates = @{
    "H20" = 2060.8235
    "I20" = 2311.5833
    "J20" = 1459
    "K20" = 2311.5833
    "L20" = 1459
    "M20" = -2064.5833
    "N20" = -1953
    "H81" = 18699.125
    "J81" = 18699.125
    "L81" = 18699.125
    "N81" = -20821.125
    "H84" = 18699.125
    "J84" = 18699.125
    "L84" = 56097.375
    "N84" = -66705.375
    "H105" = 2780823.2
    "I105" = 3750.9
    "K105" = 3750.9
    "M105" = -2003.9
    "H110" = 44831.25
    "J110" = 44831.25
    "L110" = 44831.25
    "N110" = -53011.25
    "H135" = 34817.375
    "J135" = 34817.375
    "L135" = 34817.375
    "N135" = -44957.375
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("CRP")

$updates = @{
    "H86" = 12636.786
    "I86" = 3189.25
    "J86" = 16415.8
    "K86" = 3189.25
    "L86" = 16415.8
    "M86" = -2066.25
    "N86" = -18661.8
    "H89" = 12636.786
    "I89" = 3189.25
    "J89" = 16415.8
    "K89" = 15946.25
    "L89" = 82079
    "M89" = -10330.25
    "N89" = -93311
    "H107" = 887.625
    "I107" = 752.1
    "J107" = 1113.5
    "K107" = 752.1
    "L107" = 1113.5
    "M107" = 1167.9
    "N107" = -4953.5
    "H122" = 2200.3
    "I122" = 2200.3
    "K122" = 6600.900000000001
    "M122" = -4150.900000000001
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("CUL")

$updates = @{
    "H5" = 1538.3334
    "I5" = 1405.4546
    "J5" = 3000
    "K5" = 4216.3638
    "L5" = 9000
    "M5" = -4104.3638
    "N5" = -9224
    "H12" = 108.4
    "J12" = 114.888885
    "L12" = 344.666655
    "N12" = -690.666655
    "H122" = 795.6
    "I122" = 500
    "J122" = 869.5
    "K122" = 4500
    "L122" = 7825.5
    "M122" = -2050
    "N122" = -12725.5
    "H131" = 797.62
    "I131" = 586.3333
    "J131" = 804.15466
    "K131" = 1758.9999
    "L131" = 2412.46398
    "M131" = 3281.0001
    "N131" = -12492.46398
    "H135" = 1538.3334
    "I135" = 1405.4546
    "J135" = 3000
    "K135" = 12649.0914
    "L135" = 27000
    "M135" = -10114.0914
    "N135" = -32070
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("GSM")

$updates = @{
    "H70" = 6958778.5
    "I70" = 49504
    "J70" = 8932857
    "K70" = 49504
    "L70" = 8932857
    "M70" = -49234
    "N70" = -8933397
    "H73" = 6958778.5
    "I73" = 49504
    "J73" = 8932857
    "K73" = 49504
    "L73" = 8932857
    "M73" = -48568
    "N73" = -8934729
    "H97" = 3159.375
    "I97" = 879.1667
    "J97" = 10000
    "K97" = 879.1667
    "L97" = 10000
    "M97" = -383.1667
    "N97" = -10992
    "H102" = 1930.5769
    "I102" = 1943.2609
    "K102" = 1943.2609
    "M102" = -321.2609
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("LTW")

$updates = @{
    "H22" = 2945.4546
    "I22" = 2075.125
    "J22" = 5266.3335
    "K22" = 2075.125
    "L22" = 5266.3335
    "M22" = -1780.125
    "N22" = -5856.3335
    "H27" = 2945.4546
    "I27" = 2075.125
    "J27" = 5266.3335
    "K27" = 2075.125
    "L27" = 5266.3335
    "M27" = -1968.125
    "N27" = -5480.3335
    "H46" = 2875
    "I46" = 2500
    "K46" = 2500
    "M46" = -2312
    "H100" = 3989.7
    "I100" = 1933.3334
    "J100" = 4871
    "K100" = 1933.3334
    "L100" = 4871
    "M100" = -1392.3334
    "N100" = -5953
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("WVR")

$updates = @{
    "H8" = 0
    "I8" = 0
    "K8" = 0
    "H101" = 17200.334
    "J101" = 17200.334
    "L101" = 17200.334
    "N101" = -23690.334
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

foreach ($addr in @("M8")) {
    $ws.Range($addr).ClearContents()
}

Write-Output "Updated leve profit figures across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR."